$wb = $excel.ActiveWorkbook

# Column letters A.. S mapped to indices 1..19 for Cells.Item(row, col)
$cols = 1..19

function Set-RowValues {
    param($ws, $rowNum, $values)
    for ($i = 0; $i -lt $values.Length; $i++) {
        $val = $values[$i]
        if ($null -ne $val) {
            $ws.Cells.Item($rowNum, $i + 1).Value2 = $val
        }
    }
}

# ---- Sheet: Team Order ----
$ws = $wb.Worksheets.Item("Team Order")

# Extend the bold/centered/bordered header style (currently on M1, the
# last "WeekN" header) across the new N1:S1 header cells so the new
# "Week13".."Week17"/"Total" labels pick up the same formatting.
$ws.Range("M1").Copy($ws.Range("N1:S1"))

Set-RowValues $ws 1 @($null,"Week1","Week2","Week3","Week4","Week5","Week6","Week7","Week8","Week9","Week10","Week11","Week12","Week13","Week14","Week15","Week16","Week17","Total")
Set-RowValues $ws 2 @("miles",24,37,31,33,32,34,34,38,30,30,32,34,33,34,27,34,5,522)
Set-RowValues $ws 3 @("kenneth",30,28,26,26,37,34,33,35,28,30,23,27,27,31,22,32,3,472)
Set-RowValues $ws 4 @("lamin",25,35,32,28,29,34,32,34,32,30,30,29,24,30,26,31,7,488)
Set-RowValues $ws 5 @("ned",24,32,34,33,26,33,32,30,28,28,29,32,29,33,26,26,6,481)
Set-RowValues $ws 6 @("rico",21,29,29,28,30,36,28,30,34,29,32,26,32,32,30,26,6,478)
Set-RowValues $ws 7 @("spencer",21,31,26,33,32,33,31,35,31,30,29,27,34,39,33,28,4,497)
Set-RowValues $ws 8 @("narayan",26,31,31,34,34,34,38,33,34,28,32,33,33,33,26,38,4,522)
Set-RowValues $ws 9 @("yinka",24,29,30,27,32,34,30,35,29,31,24,35,31,35,28,28,6,488)
Set-RowValues $ws 10 @("ryan",20,27,30,29,33,25,27,30,27,21,18,21,21,25,16,15,2,387)
Set-RowValues $ws 11 @("zach",26,33,31,30,25,31,21,26,29,27,27,28,24,28,26,31,6,449)
Set-RowValues $ws 12 @("abass",25,31,31,28,25,32,23,26,27,25,23,28,32,37,31,29,6,459)
Set-RowValues $ws 13 @("jordan",26,35,31,27,32,35,28,30,33,29,32,32,29,29,30,34,7,499)
Set-RowValues $ws 14 @("jordanc",27,33,29,25,33,28,32,29,26,21,25,31,25,22,17,26,4,433)
Set-RowValues $ws 15 @("hurley",28,28,24,30,31,35,31,28,32,22,25,26,24,25,27,26,4,446)

# ---- Sheet: Total Order ----
$ws = $wb.Worksheets.Item("Total Order")

# Extend the bold/centered/bordered header style (currently on M1, the
# last "WeekN" header) across the new N1:S1 header cells so the new
# "Week13".."Week17"/"Total" labels pick up the same formatting.
$ws.Range("M1").Copy($ws.Range("N1:S1"))

Set-RowValues $ws 1 @($null,"Week1","Week2","Week3","Week4","Week5","Week6","Week7","Week8","Week9","Week10","Week11","Week12","Week13","Week14","Week15","Week16","Week17","Total")
Set-RowValues $ws 2 @("ryan",20,27,30,29,33,25,27,30,27,21,18,21,21,25,16,15,2,387)
Set-RowValues $ws 3 @("jordanc",27,33,29,25,33,28,32,29,26,21,25,31,25,22,17,26,4,433)
Set-RowValues $ws 4 @("hurley",28,28,24,30,31,35,31,28,32,22,25,26,24,25,27,26,4,446)
Set-RowValues $ws 5 @("zach",26,33,31,30,25,31,21,26,29,27,27,28,24,28,26,31,6,449)
Set-RowValues $ws 6 @("abass",25,31,31,28,25,32,23,26,27,25,23,28,32,37,31,29,6,459)
Set-RowValues $ws 7 @("kenneth",30,28,26,26,37,34,33,35,28,30,23,27,27,31,22,32,3,472)
Set-RowValues $ws 8 @("rico",21,29,29,28,30,36,28,30,34,29,32,26,32,32,30,26,6,478)
Set-RowValues $ws 9 @("ned",24,32,34,33,26,33,32,30,28,28,29,32,29,33,26,26,6,481)
Set-RowValues $ws 10 @("lamin",25,35,32,28,29,34,32,34,32,30,30,29,24,30,26,31,7,488)
Set-RowValues $ws 11 @("yinka",24,29,30,27,32,34,30,35,29,31,24,35,31,35,28,28,6,488)
Set-RowValues $ws 12 @("spencer",21,31,26,33,32,33,31,35,31,30,29,27,34,39,33,28,4,497)
Set-RowValues $ws 13 @("jordan",26,35,31,27,32,35,28,30,33,29,32,32,29,29,30,34,7,499)
Set-RowValues $ws 14 @("miles",24,37,31,33,32,34,34,38,30,30,32,34,33,34,27,34,5,522)
Set-RowValues $ws 15 @("narayan",26,31,31,34,34,34,38,33,34,28,32,33,33,33,26,38,4,522)

